# Arabic Input field and one field personal address added
# Adds a new "PersonalAddress" column (Y) to the EmpInsider and
# NonEmpInsider mass-upload template sheets.

$wb = $excel.ActiveWorkbook

$empInsider = $wb.Worksheets.Item("EmpInsider")
$nonEmpInsider = $wb.Worksheets.Item("NonEmpInsider")

# --- EmpInsider (sheet1): new column Y "PersonalAddress" -------------------
# Clone the header style from the last existing header cell (X1) so the new
# header cell gets the same bold / fill formatting (style index 6).
$empInsider.Range("X1").Copy()
$empInsider.Range("Y1").PasteSpecial(-4122)
$empInsider.Range("Y1").Value = "PersonalAddress"
$empInsider.Columns.Item(25).ColumnWidth = 15.14

# --- NonEmpInsider (sheet2): new column Y "PersonalAddress" ----------------
$nonEmpInsider.Range("X1").Copy()
$nonEmpInsider.Range("Y1").PasteSpecial(-4122)
$nonEmpInsider.Range("Y1").Value = "PersonalAddress"

# --- Selections -------------------------------------------------------------
# Both sheets now show an entire-column selection on column I; select the
# non-active sheet first so the originally active sheet (EmpInsider) ends up
# active/selected again afterwards.
$nonEmpInsider.Range("I1:I1048576").Select()
$empInsider.Range("I1:I1048576").Select()
